# Update the Five Year Plan revenue_summary workbook from FY2022-FY2026 to FY2023-FY2027
# - Bumps all "FY YYYY" / "FY YYYY - FY YYYY" / "The Five Year Plan: FY YYYY - FY YYYY" labels by one year
# - Updates the plan publish date on the two summary sheets
# - Overwrites the raw data tables ("Revenue Data" and "Tax Base Data") with the refreshed source values;
#   every other sheet is formula-driven off of those two tables and recalculates automatically.

$wb = $excel.ActiveWorkbook

function BumpYears($s) {
    $result = ""
    $i = 0
    $len = $s.Length
    while ($i -lt $len) {
        if ($i + 4 -le $len) {
            $sub = $s.Substring($i, 4)
            if ($sub -match "^\d\d\d\d$") {
                $yr = [int]$sub
                $result = $result + [string]($yr + 1)
                $i = $i + 4
                continue
            }
        }
        $result = $result + $s.Substring($i,1)
        $i = $i + 1
    }
    return $result
}

function BumpCells($ws, $cells) {
    foreach ($c in $cells) {
        $old = $ws.Range($c).Value2
        $ws.Range($c).Value = BumpYears($old)
    }
}

# --- Revenue by FY ---
$ws = $wb.Worksheets.Item("Revenue by FY")
BumpCells $ws @("A2", "A6", "A19", "A32", "A45", "A58", "A71")
$ws.Range("A4").Value = "07/08/2022"

# --- Revenue by Tax ---
$ws = $wb.Worksheets.Item("Revenue by Tax")
BumpCells $ws @("A2", "A7", "A8", "A9", "A10", "A11", "A15", "A16", "A17", "A18", "A19", "A23", "A24", "A25", "A26", "A27", "A31", "A32", "A33", "A34", "A35", "A39", "A40", "A41", "A42", "A43", "A47", "A48", "A49", "A50", "A51", "A55", "A56", "A57", "A58", "A59", "A63", "A64", "A65", "A66", "A67", "A71", "A72", "A73", "A74", "A75", "A79", "A80", "A81", "A82", "A83")
$ws.Range("A4").Value = "07/08/2022"

# --- Growth by Tax ---
$ws = $wb.Worksheets.Item("Growth by Tax")
BumpCells $ws @("B2", "B9", "B10", "B11", "B12", "B13", "B16", "B17", "B18", "B19", "B20", "B23", "B24", "B25", "B26", "B27", "B30", "B31", "B32", "B33", "B34", "B37", "B38", "B39", "B40", "B41", "B44", "B45", "B46", "B47", "B48", "B51", "B52", "B53", "B54", "B55", "B58", "B59", "B60", "B61", "B62")

# --- Growth by Year ---
$ws = $wb.Worksheets.Item("Growth by Year")
BumpCells $ws @("B2", "B7", "B17", "B27", "B37", "B47")

# --- Revenue Data (raw source values; fiscal_year, Five Year Plan $, Controller $) ---
$ws = $wb.Worksheets.Item("Revenue Data")
$revenueData = @(
    @("A2", 2022),
    @("B2", 1538713000),
    @("C2", 1538713000),
    @("A3", 2023),
    @("B3", 1639878000),
    @("C3", 1614849683.236435),
    @("A4", 2024),
    @("B4", 1711213000),
    @("C4", 1675102929.591333),
    @("A5", 2025),
    @("B5", 1786677000),
    @("C5", 1748979471.934149),
    @("A6", 2026),
    @("B6", 1860646000),
    @("C6", 1821577646.923966),
    @("A7", 2027),
    @("B7", 1936560000),
    @("C7", 1896136850.325276),
    @("A8", 2022),
    @("B8", 389046376.5480998),
    @("C8", 389046376.5480998),
    @("A9", 2023),
    @("B9", 398266776),
    @("C9", 403417433.0675004),
    @("A10", 2024),
    @("B10", 409816512),
    @("C10", 413625204.1699599),
    @("A11", 2025),
    @("B11", 422766714),
    @("C11", 422199599.6194052),
    @("A12", 2026),
    @("B12", 436422079),
    @("C12", 430986990.7740095),
    @("A13", 2027),
    @("B13", 449776594),
    @("C13", 439579454.5670376),
    @("A14", 2022),
    @("B14", 634257000),
    @("C14", 634257000),
    @("A15", 2023),
    @("B15", 631530000),
    @("C15", 638993318.6292732),
    @("A16", 2024),
    @("B16", 647974000),
    @("C16", 656815497.0385745),
    @("A17", 2025),
    @("B17", 655712000),
    @("C17", 674610205.5962369),
    @("A18", 2026),
    @("B18", 677384000),
    @("C18", 692404740.9302766),
    @("A19", 2027),
    @("B19", 701331000),
    @("C19", 710199275.8517219),
    @("A20", 2022),
    @("B20", 572420007.3721925),
    @("C20", 572420007.3721925),
    @("A21", 2023),
    @("B21", 418307000),
    @("C21", 384203801.9803227),
    @("A22", 2024),
    @("B22", 418642000),
    @("C22", 377029570.7340896),
    @("A23", 2025),
    @("B23", 421782000),
    @("C23", 341259055.6704836),
    @("A24", 2026),
    @("B24", 427054000),
    @("C24", 337827253.7509435),
    @("A25", 2027),
    @("B25", 434143000),
    @("C25", 339946831.6794177),
    @("A26", 2022),
    @("B26", 88284000.00000004),
    @("C26", 88284000.00000004),
    @("A27", 2023),
    @("B27", 93140000),
    @("C27", 95261220.00572196),
    @("A28", 2024),
    @("B28", 95804000),
    @("C28", 100227699.6389275),
    @("A29", 2025),
    @("B29", 98716000),
    @("C29", 103687342.9624242),
    @("A30", 2026),
    @("B30", 101579000),
    @("C30", 106829208.7861063),
    @("A31", 2027),
    @("B31", 104403000),
    @("C31", 110091568.9779931),
    @("A32", 2022),
    @("B32", 21828000),
    @("C32", 21828000),
    @("A33", 2023),
    @("B33", 22701000),
    @("C33", 22264752.23760852),
    @("A34", 2024),
    @("B34", 23448000),
    @("C34", 23739332.37084913),
    @("A35", 2025),
    @("B35", 24255000),
    @("C35", 25002694.18771732),
    @("A36", 2026),
    @("B36", 25026000),
    @("C36", 26270657.81741631),
    @("A37", 2027),
    @("B37", 25782000),
    @("C37", 27586933.45817689),
    @("A38", 2022),
    @("B38", 25658000),
    @("C38", 25658000),
    @("A39", 2023),
    @("B39", 31726000),
    @("C39", 33611983.67133909),
    @("A40", 2024),
    @("B40", 33379000),
    @("C40", 32024585.43860196),
    @("A41", 2025),
    @("B41", 35916000),
    @("C41", 33441747.90452376),
    @("A42", 2026),
    @("B42", 38171000),
    @("C42", 34514547.34986574),
    @("A43", 2027),
    @("B43", 40133000),
    @("C43", 35156762.74384565)
)
foreach ($item in $revenueData) {
    $ws.Range($item[0]).Value = $item[1]
}

# --- Tax Base Data (raw source values; fiscal_year end date serial + tax base $ by type) ---
$ws = $wb.Worksheets.Item("Tax Base Data")
$taxBaseData = @(
    @("A2", 44562),
    @("B2", 56182260182.24119),
    @("C2", 19452318827.40499),
    @("D2", 19570295066.11193),
    @("E2", 12102430388.69258),
    @("F2", 7467864677.419355),
    @("G2", 17462477345.09434),
    @("H2", 392373333.3333335),
    @("I2", 436560000),
    @("J2", 954189956.7658008),
    @("A3", 44927),
    @("B3", 59720772309.03976),
    @("C3", 20170871653.37502),
    @("D3", 20260831716.41102),
    @("E3", 12560229534.07098),
    @("F3", 7700602182.340047),
    @("G3", 11720677302.63339),
    @("H3", 423383200.0254309),
    @("I3", 445295044.7521703),
    @("J3", 1263038616.839737),
    @("A4", 45292),
    @("B4", 61949072839.91616),
    @("C4", 20681260208.498),
    @("D4", 20858335293.48096),
    @("E4", 12952978951.71442),
    @("F4", 7905356341.766534),
    @("G4", 11501817289.02042),
    @("H4", 445456442.8396778),
    @("I4", 474786647.4169827),
    @("J4", 1203388901.195023),
    @("A5", 45658),
    @("B5", 64681193488.68895),
    @("C5", 21109979980.97026),
    @("D5", 21455054496.46941),
    @("E5", 13345301830.43236),
    @("F5", 8109752666.037045),
    @("G5", 10410587421.30823),
    @("H5", 460832635.3885521),
    @("I5", 500053883.7543464),
    @("J5", 1256641661.826385),
    @("A6", 46023),
    @("B6", 67366037238.31236),
    @("C6", 21549349538.70048),
    @("D6", 22051770174.39204),
    @("E6", 13737623880.12693),
    @("F6", 8314146294.265119),
    @("G6", 10305895477.45404),
    @("H6", 474796483.4938059),
    @("I6", 525413156.3483262),
    @("J6", 1296954281.897856),
    @("A7", 46388),
    @("B7", 70123404228.00577),
    @("C7", 21978972728.35188),
    @("D7", 22648485841.54984),
    @("E7", 14129945924.74565),
    @("F7", 8518539916.804188),
    @("G7", 10370556183.02067),
    @("H7", 489295862.1244138),
    @("I7", 551738669.1635377),
    @("J7", 1321086830.897552)
)
foreach ($item in $taxBaseData) {
    $ws.Range($item[0]).Value = $item[1]
}
